$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 7).Value = 4.408429333333333
$ws.Cells.Item(2, 8).Value = 13.225288
$ws.Cells.Item(2, 9).Value = 0.05980478019486075
$ws.Cells.Item(2, 10).Value = 0.05980478019486075
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 13).Value = 19.68532233333334
$ws.Cells.Item(2, 14).Value = 59.055967
$ws.Cells.Item(2, 15).Value = 0.5039562241205585
$ws.Cells.Item(2, 16).Value = 0.5039562241205585
$ws.Cells.Item(2, 17).Value = 86.78135241038845
$ws.Cells.Item(2, 18).Value = 781.032171693496
$ws.Cells.Item(2, 19).Value = 0.03013899121136198
$ws.Cells.Item(2, 20).Value = 0.03013899121136198
$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 7).Value = 4.408429333333333
$ws.Cells.Item(3, 8).Value = 13.225288
$ws.Cells.Item(3, 9).Value = 0.05980478019486075
$ws.Cells.Item(3, 10).Value = 0.05980478019486075
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 13).Value = 1.484487666666667
$ws.Cells.Item(3, 14).Value = 4.453463
$ws.Cells.Item(3, 15).Value = 0.03800378711503639
$ws.Cells.Item(3, 16).Value = 0.0380037871150364
$ws.Cells.Item(3, 17).Value = 6.544258974704888
$ws.Cells.Item(3, 18).Value = 58.898330772344
$ws.Cells.Item(3, 19).Value = 0.002272808134987033
$ws.Cells.Item(3, 20).Value = 0.002272808134987033
$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 7).Value = 4.408429333333333
$ws.Cells.Item(4, 8).Value = 13.225288
$ws.Cells.Item(4, 9).Value = 0.05980478019486075
$ws.Cells.Item(4, 10).Value = 0.05980478019486075
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 13).Value = 17.76686133333333
$ws.Cells.Item(4, 14).Value = 53.300584
$ws.Cells.Item(4, 15).Value = 0.4548424557345857
$ws.Cells.Item(4, 16).Value = 0.4548424557345858
$ws.Cells.Item(4, 17).Value = 78.32395266313245
$ws.Cells.Item(4, 18).Value = 704.915573968192
$ws.Cells.Item(4, 19).Value = 0.02720175308849758
$ws.Cells.Item(4, 20).Value = 0.02720175308849758
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 7).Value = 4.408429333333333
$ws.Cells.Item(5, 8).Value = 13.225288
$ws.Cells.Item(5, 9).Value = 0.05980478019486075
$ws.Cells.Item(5, 10).Value = 0.05980478019486075
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 13).Value = 0.1249006666666667
$ws.Cells.Item(5, 14).Value = 0.374702
$ws.Cells.Item(5, 15).Value = 0.003197533029819349
$ws.Cells.Item(5, 16).Value = 0.003197533029819349
$ws.Cells.Item(5, 17).Value = 0.5506157626862221
$ws.Cells.Item(5, 18).Value = 4.955541864175999
$ws.Cells.Item(5, 19).Value = 0.0001912277600141533
$ws.Cells.Item(5, 20).Value = 0.0001912277600141533
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 7).Value = 56.15338233333333
$ws.Cells.Item(6, 8).Value = 168.460147
$ws.Cells.Item(6, 9).Value = 0.7617771395926449
$ws.Cells.Item(6, 10).Value = 0.7617771395926448
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 13).Value = 19.68532233333334
$ws.Cells.Item(6, 14).Value = 59.055967
$ws.Cells.Item(6, 15).Value = 0.5039562241205585
$ws.Cells.Item(6, 16).Value = 0.5039562241205585
$ws.Cells.Item(6, 17).Value = 1105.397431338572
$ws.Cells.Item(6, 18).Value = 9948.576882047149
$ws.Cells.Item(6, 19).Value = 0.3839023308904689
$ws.Cells.Item(6, 20).Value = 0.3839023308904689
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 7).Value = 56.15338233333333
$ws.Cells.Item(7, 8).Value = 168.460147
$ws.Cells.Item(7, 9).Value = 0.7617771395926449
$ws.Cells.Item(7, 10).Value = 0.7617771395926448
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 13).Value = 1.484487666666667
$ws.Cells.Item(7, 14).Value = 4.453463
$ws.Cells.Item(7, 15).Value = 0.03800378711503639
$ws.Cells.Item(7, 16).Value = 0.0380037871150364
$ws.Cells.Item(7, 17).Value = 83.35900351545122
$ws.Cells.Item(7, 18).Value = 750.231031639061
$ws.Cells.Item(7, 19).Value = 0.02895041624218024
$ws.Cells.Item(7, 20).Value = 0.02895041624218024
$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 7).Value = 56.15338233333333
$ws.Cells.Item(8, 8).Value = 168.460147
$ws.Cells.Item(8, 9).Value = 0.7617771395926449
$ws.Cells.Item(8, 10).Value = 0.7617771395926448
$ws.Cells.Item(8, 11).Value = 3
$ws.Cells.Item(8, 13).Value = 17.76686133333333
$ws.Cells.Item(8, 14).Value = 53.300584
$ws.Cells.Item(8, 15).Value = 0.4548424557345857
$ws.Cells.Item(8, 16).Value = 0.4548424557345858
$ws.Cells.Item(8, 17).Value = 997.6693573139831
$ws.Cells.Item(8, 18).Value = 8979.024215825848
$ws.Cells.Item(8, 19).Value = 0.3464885848947869
$ws.Cells.Item(8, 20).Value = 0.3464885848947869
$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 7).Value = 56.15338233333333
$ws.Cells.Item(9, 8).Value = 168.460147
$ws.Cells.Item(9, 9).Value = 0.7617771395926449
$ws.Cells.Item(9, 10).Value = 0.7617771395926448
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 13).Value = 0.1249006666666667
$ws.Cells.Item(9, 14).Value = 0.374702
$ws.Cells.Item(9, 15).Value = 0.003197533029819349
$ws.Cells.Item(9, 16).Value = 0.003197533029819349
$ws.Cells.Item(9, 17).Value = 7.013594889021555
$ws.Cells.Item(9, 18).Value = 63.122354001194
$ws.Cells.Item(9, 19).Value = 0.002435807565208787
$ws.Cells.Item(9, 20).Value = 0.002435807565208787
$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 7).Value = 2.790021666666667
$ws.Cells.Item(10, 8).Value = 8.370065
$ws.Cells.Item(10, 9).Value = 0.03784945156141002
$ws.Cells.Item(10, 10).Value = 0.03784945156141001
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 13).Value = 19.68532233333334
$ws.Cells.Item(10, 14).Value = 59.055967
$ws.Cells.Item(10, 15).Value = 0.5039562241205585
$ws.Cells.Item(10, 16).Value = 0.5039562241205585
$ws.Cells.Item(10, 17).Value = 54.92247582531723
$ws.Cells.Item(10, 18).Value = 494.3022824278551
$ws.Cells.Item(10, 19).Value = 0.01907446669392217
$ws.Cells.Item(10, 20).Value = 0.01907446669392217
$ws.Cells.Item(11, 5).Value = 3
$ws.Cells.Item(11, 7).Value = 2.790021666666667
$ws.Cells.Item(11, 8).Value = 8.370065
$ws.Cells.Item(11, 9).Value = 0.03784945156141002
$ws.Cells.Item(11, 10).Value = 0.03784945156141001
$ws.Cells.Item(11, 11).Value = 3
$ws.Cells.Item(11, 13).Value = 1.484487666666667
$ws.Cells.Item(11, 14).Value = 4.453463
$ws.Cells.Item(11, 15).Value = 0.03800378711503639
$ws.Cells.Item(11, 16).Value = 0.0380037871150364
$ws.Cells.Item(11, 17).Value = 4.141752753899445
$ws.Cells.Item(11, 18).Value = 37.275774785095
$ws.Cells.Item(11, 19).Value = 0.001438422499560708
$ws.Cells.Item(11, 20).Value = 0.001438422499560708
$ws.Cells.Item(12, 5).Value = 3
$ws.Cells.Item(12, 7).Value = 2.790021666666667
$ws.Cells.Item(12, 8).Value = 8.370065
$ws.Cells.Item(12, 9).Value = 0.03784945156141002
$ws.Cells.Item(12, 10).Value = 0.03784945156141001
$ws.Cells.Item(12, 11).Value = 3
$ws.Cells.Item(12, 13).Value = 17.76686133333333
$ws.Cells.Item(12, 14).Value = 53.300584
$ws.Cells.Item(12, 15).Value = 0.4548424557345857
$ws.Cells.Item(12, 16).Value = 0.4548424557345858
$ws.Cells.Item(12, 17).Value = 49.56992806866223
$ws.Cells.Item(12, 18).Value = 446.12935261796
$ws.Cells.Item(12, 19).Value = 0.01721553749639898
$ws.Cells.Item(12, 20).Value = 0.01721553749639898
$ws.Cells.Item(13, 5).Value = 3
$ws.Cells.Item(13, 7).Value = 2.790021666666667
$ws.Cells.Item(13, 8).Value = 8.370065
$ws.Cells.Item(13, 9).Value = 0.03784945156141002
$ws.Cells.Item(13, 10).Value = 0.03784945156141001
$ws.Cells.Item(13, 11).Value = 3
$ws.Cells.Item(13, 13).Value = 0.1249006666666667
$ws.Cells.Item(13, 14).Value = 0.374702
$ws.Cells.Item(13, 15).Value = 0.003197533029819349
$ws.Cells.Item(13, 16).Value = 0.003197533029819349
$ws.Cells.Item(13, 17).Value = 0.3484755661811111
$ws.Cells.Item(13, 18).Value = 3.13628009563
$ws.Cells.Item(13, 19).Value = 0.0001210248715281561
$ws.Cells.Item(13, 20).Value = 0.0001210248715281561
$ws.Cells.Item(14, 5).Value = 3
$ws.Cells.Item(14, 7).Value = 10.36182833333333
$ws.Cells.Item(14, 8).Value = 31.085485
$ws.Cells.Item(14, 9).Value = 0.1405686286510843
$ws.Cells.Item(14, 10).Value = 0.1405686286510843
$ws.Cells.Item(14, 11).Value = 3
$ws.Cells.Item(14, 13).Value = 19.68532233333334
$ws.Cells.Item(14, 14).Value = 59.055967
$ws.Cells.Item(14, 15).Value = 0.5039562241205585
$ws.Cells.Item(14, 16).Value = 0.5039562241205585
$ws.Cells.Item(14, 17).Value = 203.9759307043328
$ws.Cells.Item(14, 18).Value = 1835.783376338995
$ws.Cells.Item(14, 19).Value = 0.0708404353248054
$ws.Cells.Item(14, 20).Value = 0.07084043532480538
$ws.Cells.Item(15, 5).Value = 3
$ws.Cells.Item(15, 7).Value = 10.36182833333333
$ws.Cells.Item(15, 8).Value = 31.085485
$ws.Cells.Item(15, 9).Value = 0.1405686286510843
$ws.Cells.Item(15, 10).Value = 0.1405686286510843
$ws.Cells.Item(15, 11).Value = 3
$ws.Cells.Item(15, 13).Value = 1.484487666666667
$ws.Cells.Item(15, 14).Value = 4.453463
$ws.Cells.Item(15, 15).Value = 0.03800378711503639
$ws.Cells.Item(15, 16).Value = 0.0380037871150364
$ws.Cells.Item(15, 17).Value = 15.38200636495056
$ws.Cells.Item(15, 18).Value = 138.438057284555
$ws.Cells.Item(15, 19).Value = 0.005342140238308413
$ws.Cells.Item(15, 20).Value = 0.005342140238308413
$ws.Cells.Item(16, 5).Value = 3
$ws.Cells.Item(16, 7).Value = 10.36182833333333
$ws.Cells.Item(16, 8).Value = 31.085485
$ws.Cells.Item(16, 9).Value = 0.1405686286510843
$ws.Cells.Item(16, 10).Value = 0.1405686286510843
$ws.Cells.Item(16, 11).Value = 3
$ws.Cells.Item(16, 13).Value = 17.76686133333333
$ws.Cells.Item(16, 14).Value = 53.300584
$ws.Cells.Item(16, 15).Value = 0.4548424557345857
$ws.Cells.Item(16, 16).Value = 0.4548424557345858
$ws.Cells.Item(16, 17).Value = 184.0971671581378
$ws.Cells.Item(16, 18).Value = 1656.87450442324
$ws.Cells.Item(16, 19).Value = 0.06393658025490223
$ws.Cells.Item(16, 20).Value = 0.06393658025490222
$ws.Cells.Item(17, 5).Value = 3
$ws.Cells.Item(17, 7).Value = 10.36182833333333
$ws.Cells.Item(17, 8).Value = 31.085485
$ws.Cells.Item(17, 9).Value = 0.1405686286510843
$ws.Cells.Item(17, 10).Value = 0.1405686286510843
$ws.Cells.Item(17, 11).Value = 3
$ws.Cells.Item(17, 13).Value = 0.1249006666666667
$ws.Cells.Item(17, 14).Value = 0.374702
$ws.Cells.Item(17, 15).Value = 0.003197533029819349
$ws.Cells.Item(17, 16).Value = 0.003197533029819349
$ws.Cells.Item(17, 17).Value = 1.294199266718889
$ws.Cells.Item(17, 18).Value = 11.64779340047
$ws.Cells.Item(17, 19).Value = 0.0004494728330682525
$ws.Cells.Item(17, 20).Value = 0.0004494728330682525
